$d = $word.ActiveDocument

# Locate the empty paragraph that immediately precedes the "Le déroulement
# du projet ..." paragraph (identified by its distinctive leading text),
# then append a brand-new paragraph right after it.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "Le déroulement du projet*") {
        $target = $d.Paragraphs.Item($i - 1)
        break
    }
}

$newPara = $target.Range.InsertParagraphAfter()

# Re-fetch the freshly inserted (now empty) paragraph that sits right
# after $target and right before the "Le déroulement ..." paragraph.
$insertedIndex = $target.Index + 1
$inserted = $d.Paragraphs.Item($insertedIndex)
$insertedRange = $inserted.Range

$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:rtl w:val="0"/></w:rPr><w:tab/><w:t xml:space="preserve">Pour le projet, si nous avions eu plus de temps, nous aurions essayé de gérer les égalités dans le jeu. Puisque pour le moment, si deux joueurs rentrent le même nombre de coups, seul le 1er sera choisi pour jouer. De même pour les scores finaux, seul le 1er joueur qui sera à égalité sera désigné comme gagnant.</w:t></w:r></w:p>'
$insertedRange.InsertXML($xml)

# InsertXML silently drops explicit-zero indent attributes during
# re-serialization, and drops the paragraph-mark run formatting; restore
# both so the paragraph properties match the original authoring exactly.
$inserted2 = $d.Paragraphs.Item($insertedIndex)
$inserted2.Range.ParagraphFormat.LeftIndent = 0
$inserted2.Range.ParagraphFormat.FirstLineIndent = 0
$inserted2.Range.Font.Size = 12
$inserted2.Range.Font.SizeBi = 12

Write-Host "Inserted paragraph at index" $insertedIndex "- text:" $inserted2.Range.Text
